$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the top H1
#    ("Play A Night of Mystery for Free") paragraph.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
$titleRange.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"
$metaRange = $metaPara.Range

$metaXmlFragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of A Night of Mystery, a traditional online slot game with engaging atmosphere and interesting winning potential. Play it for free!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$metaRange.InsertXML($metaXmlFragment)

# ---------------------------------------------------------------------------
# 2) Drop the trailing duplicate "Play A Night of Mystery for Free" (bold)
#    paragraph and turn the final italic paragraph's text into the AI
#    image-generation prompt, keeping its italic formatting intact.
# ---------------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs($lastIndex - 1)
$dupTitlePara.Range.Delete()

$promptPara = $d.Paragraphs($d.Paragraphs.Count)
$promptRange = $promptPara.Range
$replaceRange = $d.Range($promptRange.Start, $promptRange.End - 1)

$promptText = 'Prompt: "A Night of Mystery" Feature Image Can you create an eye-catching feature image for "A Night of Mystery"? The image should be in a cartoon style and should feature a happy Maya warrior with glasses as it is the main character in the game. Please make sure the image is bright and bold to attract potential players. You can showcase the theme of unsolved murders and feature the burlesque dancers as well. Use your creativity to ensure that the image stands out and reflects the fun and exciting gaming experience of "A Night of Mystery".'

$replaceRange.Text = $promptText
